$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("additive_200epochs", "ranking_test2", 3, 0.001,  "mean", "product", 198.448, 0.767, 12.781, 0.237, 17.706),
    @("additive_200epochs", "ranking_test2", 3, 0.0001, "mean", "product", 287.345, 0.613, 16.588, 0.175, 19.482),
    @("additive_200epochs", "ranking_test2", 5, 0.0001, "mean", "product", 167.498, 0.773, 12.605, 0.337, 17.072),
    @("additive_200epochs", "ranking_test2", 2, 0.0001, "sum",  "product", 281.699, 0.619, 16.462, 0.248, 18.322),
    @("additive_200epochs", "ranking_test2", 3, 0.001,  "max",  "mean",    223.099, 0.712, 14.326, 0.273, 17.479),
    @("additive_200epochs", "ranking_test2", 3, 0.0001, "max",  "mean",    355.712, 0.507, 18.72,  0.187, 19.182)
)

$startRow = 91
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
